# Refresh the cryptos worksheet: for each coin row, update the Price (D) and
# Volume(1h) (E) columns with the latest scraped figures, and for the handful of
# rows whose rank order changed, swap the Coin (B) / Link (C) / Price (D) /
# Volume(1h) (E) values so each row again reflects the correct coin.
#
# Every value is written with a leading apostrophe so Excel stores it as literal
# text (quote-prefixed), matching how this sheet already keeps its Price/Volume
# columns as strings -- without it, values like "2.30" would be re-interpreted
# as the number 2.3 and lose their trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'51.074.37"
$ws.Range("E2").Value = "'  -0.95%  "

# Row 3
$ws.Range("D3").Value = "'2.939.19"
$ws.Range("E3").Value = "'  -1.49%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'373.81"
$ws.Range("E5").Value = "'  -1.59%  "

# Row 6
$ws.Range("D6").Value = "'100.71"
$ws.Range("E6").Value = "'  -3.91%  "

# Row 7
$ws.Range("D7").Value = "'0.535"
$ws.Range("E7").Value = "'  -1.00%  "

# Row 8
$ws.Range("E8").Value = "'  +0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "'  -1.88%  "

# Row 10
$ws.Range("D10").Value = "'36.17"
$ws.Range("E10").Value = "'  -2.85%  "

# Row 11
$ws.Range("D11").Value = "'0.139"
$ws.Range("E11").Value = "'  -0.74%  "

# Row 12
$ws.Range("D12").Value = "'0.0848"
$ws.Range("E12").Value = "'  +0.33%  "

# Row 13
$ws.Range("D13").Value = "'3.402.71"
$ws.Range("E13").Value = "'  -1.29%  "

# Row 14
$ws.Range("D14").Value = "'17.97"
$ws.Range("E14").Value = "'  -2.49%  "

# Row 15
$ws.Range("D15").Value = "'7.52"
$ws.Range("E15").Value = "'  -0.97%  "

# Row 16
$ws.Range("D16").Value = "'2.943.07"
$ws.Range("E16").Value = "'  -1.39%  "

# Row 17
$ws.Range("D17").Value = "'11.12"
$ws.Range("E17").Value = "'  +49.45%  "

# Row 18
$ws.Range("D18").Value = "'0.974"
$ws.Range("E18").Value = "'  -0.10%  "

# Row 19
$ws.Range("D19").Value = "'50.993.80"
$ws.Range("E19").Value = "'  -1.00%  "

# Row 20
$ws.Range("D20").Value = "'3.15"

# Row 21
$ws.Range("D21").Value = "'12.39"
$ws.Range("E21").Value = "'  -4.36%  "

# Row 22
$ws.Range("D22").Value = "'0.0₃0957"
$ws.Range("E22").Value = "'  -0.59%  "

# Row 23
$ws.Range("B23").Value = "'BitcoinCash"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'264.67"
$ws.Range("E23").Value = "'  +1.04%  "

# Row 24
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'68.69"
$ws.Range("E24").Value = "'  -0.97%  "

# Row 25
$ws.Range("D25").Value = "'3.12"
$ws.Range("E25").Value = "'  +10.00%  "

# Row 26
$ws.Range("D26").Value = "'8.07"
$ws.Range("E26").Value = "'  -0.70%  "

# Row 27
$ws.Range("D27").Value = "'7.49"
$ws.Range("E27").Value = "'  -2.26%  "

# Row 28
$ws.Range("E28").Value = "'  -0.01%  "

# Row 29
$ws.Range("D29").Value = "'25.62"
$ws.Range("E29").Value = "'  -1.07%  "

# Row 30
$ws.Range("B30").Value = "'Kaspa"
$ws.Range("C30").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.163"
$ws.Range("E30").Value = "'  -4.11%  "

# Row 31
$ws.Range("B31").Value = "'Hedera"
$ws.Range("C31").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.111"
$ws.Range("E31").Value = "'  -5.23%  "

# Row 32
$ws.Range("D32").Value = "'9.92"
$ws.Range("E32").Value = "'  +0.46%  "

# Row 33
$ws.Range("D33").Value = "'50.66"
$ws.Range("E33").Value = "'  -0.81%  "

# Row 34
$ws.Range("E34").Value = "'  -1.87%  "

# Row 35
$ws.Range("D35").Value = "'33.22"
$ws.Range("E35").Value = "'  -5.45%  "

# Row 36
$ws.Range("D36").Value = "'0.0442"
$ws.Range("E36").Value = "'  -0.64%  "

# Row 37
$ws.Range("E37").Value = "'  -0.24%  "

# Row 38
$ws.Range("D38").Value = "'3.16"
$ws.Range("E38").Value = "'  +4.07%  "

# Row 39
$ws.Range("E39").Value = "'  -0.08%  "

# Row 40
$ws.Range("D40").Value = "'16.32"
$ws.Range("E40").Value = "'  -5.16%  "

# Row 41
$ws.Range("B41").Value = "'Stacks"
$ws.Range("C41").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.48"
$ws.Range("E41").Value = "'  -4.00%  "

# Row 42
$ws.Range("B42").Value = "'ARBITRUM"
$ws.Range("C42").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.79"
$ws.Range("E42").Value = "'  -3.50%  "

# Row 43
$ws.Range("D43").Value = "'119.87"
$ws.Range("E43").Value = "'  -4.44%  "

# Row 44
$ws.Range("D44").Value = "'21.07"
$ws.Range("E44").Value = "'  -2.93%  "

# Row 45
$ws.Range("B45").Value = "'TheGraph"
$ws.Range("C45").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.274"
$ws.Range("E45").Value = "'  -3.74%  "

# Row 46
$ws.Range("B46").Value = "'WEMIXToken"
$ws.Range("C46").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.03"
$ws.Range("E46").Value = "'  -1.34%  "

# Row 47
$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "'  +2.57%  "

# Row 48
$ws.Range("B48").Value = "'ApeXProtocol"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.30"
$ws.Range("E48").Value = "'  -3.57%  "

# Row 49
$ws.Range("B49").Value = "'Maker"
$ws.Range("C49").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'1.990.84"
$ws.Range("E49").Value = "'  -2.13%  "

# Row 50
$ws.Range("E50").Value = "'  -2.99%  "

# Row 51
$ws.Range("D51").Value = "'1.30"
$ws.Range("E51").Value = "'  +0.80%  "
